$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 156
$ws.Range("I5").Value = 162.5
$ws.Range("J5").Value = 149.5
$ws.Range("K5").Value = 162.5
$ws.Range("L5").Value = 149.5
$ws.Range("M5").Value = -47.5
$ws.Range("N5").Value = -379.5
$ws.Range("H17").Value = 100000
$ws.Range("J17").Value = 100000
$ws.Range("L17").Value = 300000
$ws.Range("N17").Value = -300336
$ws.Range("H28").Value = 1222.8462
$ws.Range("I28").Value = 1249.75
$ws.Range("K28").Value = 1249.75
$ws.Range("M28").Value = -764.75
$ws.Range("H58").Value = 1401.5
$ws.Range("I58").Value = 80
$ws.Range("J58").Value = 2723
$ws.Range("K58").Value = 240
$ws.Range("L58").Value = 8169
$ws.Range("M58").Value = -90
$ws.Range("N58").Value = -8469
$ws.Range("H70").Value = 3983.75
$ws.Range("J70").Value = 4649
$ws.Range("L70").Value = 13947
$ws.Range("N70").Value = -14487
$ws.Range("H73").Value = 3983.75
$ws.Range("J73").Value = 4649
$ws.Range("L73").Value = 13947
$ws.Range("N73").Value = -15819
$ws.Range("H74").Value = 68414.8
$ws.Range("I74").Value = 111238.5
$ws.Range("J74").Value = 13911.909
$ws.Range("K74").Value = 111238.5
$ws.Range("L74").Value = 13911.909
$ws.Range("M74").Value = -110302.5
$ws.Range("N74").Value = -15783.909
$ws.Range("H77").Value = 68414.8
$ws.Range("I77").Value = 111238.5
$ws.Range("J77").Value = 13911.909
$ws.Range("K77").Value = 556192.5
$ws.Range("L77").Value = 69559.545
$ws.Range("M77").Value = -551512.5
$ws.Range("N77").Value = -78919.545
$ws.Range("H92").Value = 590.8570999999999
$ws.Range("I92").Value = 615.1667
$ws.Range("K92").Value = 615.1667
$ws.Range("M92").Value = 632.8333
$ws.Range("H103").Value = 1116
$ws.Range("J103").Value = 1281.8334
$ws.Range("L103").Value = 3845.5002
$ws.Range("N103").Value = -5017.5002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6696.778
$ws.Range("I32").Value = 5632.56
$ws.Range("K32").Value = 5632.56
$ws.Range("M32").Value = -5345.56
$ws.Range("H61").Value = 5388.5
$ws.Range("H69").Value = 348787.5
$ws.Range("J69").Value = 348787.5
$ws.Range("L69").Value = 348787.5
$ws.Range("N69").Value = -350285.5
$ws.Range("H72").Value = 348787.5
$ws.Range("J72").Value = 348787.5
$ws.Range("L72").Value = 1046362.5
$ws.Range("N72").Value = -1053850.5
$ws.Range("H136").Value = 5388.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 603.65
$ws.Range("I80").Value = 917.36365
$ws.Range("J80").Value = 220.22223
$ws.Range("K80").Value = 917.36365
$ws.Range("L80").Value = 220.22223
$ws.Range("M80").Value = 80.63634999999999
$ws.Range("N80").Value = -2216.22223
$ws.Range("H83").Value = 603.65
$ws.Range("I83").Value = 917.36365
$ws.Range("J83").Value = 220.22223
$ws.Range("K83").Value = 4586.81825
$ws.Range("L83").Value = 1101.11115
$ws.Range("M83").Value = 405.1817499999997
$ws.Range("N83").Value = -11085.11115
$ws.Range("H86").Value = 1469.8
$ws.Range("I86").Value = 1469.8
$ws.Range("K86").Value = 1469.8
$ws.Range("M86").Value = -346.8
$ws.Range("H89").Value = 1469.8
$ws.Range("I89").Value = 1469.8
$ws.Range("K89").Value = 7349
$ws.Range("M89").Value = -1733
$ws.Range("H94").Value = 7822.222
$ws.Range("I94").Value = 7822.222
$ws.Range("K94").Value = 7822.222
$ws.Range("M94").Value = -7371.222
$ws.Range("H107").Value = 4518.6665
$ws.Range("I107").Value = 3524.75
$ws.Range("K107").Value = 3524.75
$ws.Range("M107").Value = -1604.75
$ws.Range("H134").Value = 2265.4443
$ws.Range("I134").Value = 2423.75
$ws.Range("K134").Value = 7271.25
$ws.Range("M134").Value = -4736.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 19999.5
$ws.Range("I86").Value = 19999.5
$ws.Range("K86").Value = 19999.5
$ws.Range("M86").Value = -18876.5
$ws.Range("H89").Value = 19999.5
$ws.Range("I89").Value = 19999.5
$ws.Range("K89").Value = 99997.5
$ws.Range("M89").Value = -94381.5
$ws.Range("H99").Value = 2456
$ws.Range("I99").Value = 1994
$ws.Range("K99").Value = 1994
$ws.Range("M99").Value = -496
$ws.Range("H107").Value = 998.75
$ws.Range("I107").Value = 818
$ws.Range("J107").Value = 1300
$ws.Range("K107").Value = 818
$ws.Range("L107").Value = 1300
$ws.Range("M107").Value = 1102
$ws.Range("N107").Value = -5140
$ws.Range("H126").Value = 2456
$ws.Range("I126").Value = 1994
$ws.Range("K126").Value = 5982
$ws.Range("M126").Value = -3512

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H137").Value = 3766.2222
$ws.Range("I137").Value = 2022.25
$ws.Range("J137").Value = 5161.4
$ws.Range("K137").Value = 6066.75
$ws.Range("L137").Value = 15484.2
$ws.Range("M137").Value = -966.75
$ws.Range("N137").Value = -25684.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 3500
$ws.Range("J14").Value = 3500
$ws.Range("L14").Value = 3500
$ws.Range("N14").Value = -3836
$ws.Range("H22").Value = 341
$ws.Range("J22").Value = 341
$ws.Range("L22").Value = 341
$ws.Range("N22").Value = -1399
$ws.Range("H122").Value = 4587.5713
$ws.Range("I122").Value = 2869
$ws.Range("K122").Value = 8607
$ws.Range("M122").Value = -6157

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6482.514
$ws.Range("I7").Value = 3886.0715
$ws.Range("J7").Value = 8213.477000000001
$ws.Range("K7").Value = 3886.0715
$ws.Range("L7").Value = 8213.477000000001
$ws.Range("M7").Value = -3774.0715
$ws.Range("N7").Value = -8437.477000000001
$ws.Range("H22").Value = 718.4286
$ws.Range("I22").Value = 676.6667
$ws.Range("J22").Value = 793.6
$ws.Range("K22").Value = 676.6667
$ws.Range("L22").Value = 793.6
$ws.Range("M22").Value = -381.6667
$ws.Range("N22").Value = -1383.6
$ws.Range("H27").Value = 718.4286
$ws.Range("I27").Value = 676.6667
$ws.Range("J27").Value = 793.6
$ws.Range("K27").Value = 676.6667
$ws.Range("L27").Value = 793.6
$ws.Range("M27").Value = -569.6667
$ws.Range("N27").Value = -1007.6
$ws.Range("H40").Value = 4569.6
$ws.Range("I40").Value = 4167.3335
$ws.Range("K40").Value = 4167.3335
$ws.Range("M40").Value = -4031.3335
$ws.Range("H82").Value = 2011.1333
$ws.Range("I82").Value = 1673.625
$ws.Range("J82").Value = 2396.8572
$ws.Range("K82").Value = 1673.625
$ws.Range("L82").Value = 2396.8572
$ws.Range("M82").Value = -1312.625
$ws.Range("N82").Value = -3118.8572
$ws.Range("H85").Value = 2011.1333
$ws.Range("I85").Value = 1673.625
$ws.Range("J85").Value = 2396.8572
$ws.Range("K85").Value = 1673.625
$ws.Range("L85").Value = 2396.8572
$ws.Range("M85").Value = -425.625
$ws.Range("N85").Value = -4892.8572
$ws.Range("H122").Value = 7021.4136
$ws.Range("I122").Value = 7423.778
$ws.Range("K122").Value = 22271.334
$ws.Range("M122").Value = -19821.334
$ws.Range("H126").Value = 6482.514
$ws.Range("I126").Value = 3886.0715
$ws.Range("J126").Value = 8213.477000000001
$ws.Range("K126").Value = 11658.2145
$ws.Range("L126").Value = 24640.431
$ws.Range("M126").Value = -9188.2145
$ws.Range("N126").Value = -29580.431
